$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.901.78"
$ws.Range("E2").Value = "  +7.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.760.96"
$ws.Range("E3").Value = "  +5.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.03"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9975"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3824"
$ws.Range("E7").Value = "  +2.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3619"
$ws.Range("E8").Value = "  +4.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.52"
$ws.Range("E9").Value = "  +4.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.228"
$ws.Range("E10").Value = "  +4.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07640"
$ws.Range("E11").Value = "  +5.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9990"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.72"
$ws.Range("E13").Value = "  +4.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.480"
$ws.Range("E14").Value = "  +7.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.087"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.769.01"
$ws.Range("E16").Value = "  +6.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001158"
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9970"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06832"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "87.00"
$ws.Range("E20").Value = "  +6.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.62"
$ws.Range("E21").Value = "  +6.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.511"
$ws.Range("E22").Value = "  +6.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.77"
$ws.Range("E23").Value = "  +5.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "25.855.79"
$ws.Range("E24").Value = "  +7.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.442"
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.958"
$ws.Range("E26").Value = "  +9.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.65"
$ws.Range("E27").Value = "  +5.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.93"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.966.83"
$ws.Range("E29").Value = "  +6.24%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.12"
$ws.Range("E30").Value = "  +5.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.230"
$ws.Range("E31").Value = "  +24.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.308"
$ws.Range("E32").Value = "  +15.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.302"
$ws.Range("E33").Value = "  +4.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "14.10"
$ws.Range("E34").Value = "  +14.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.822"
$ws.Range("E35").Value = "  +5.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08729"
$ws.Range("E36").Value = "  +3.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.706"
$ws.Range("E37").Value = "  +6.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06760"
$ws.Range("E38").Value = "  +5.86%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02500"
$ws.Range("E39").Value = "  +7.70%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.327"
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("E41").Value = "  +7.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.292"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6553"
$ws.Range("E43").Value = "  +7.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.35"
$ws.Range("E44").Value = "  +7.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9962"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6315"
$ws.Range("E46").Value = "  +5.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.907"
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.165"
$ws.Range("E48").Value = "  +7.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.59"
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07493"
$ws.Range("E50").Value = "  +5.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.80"
$ws.Range("E51").Value = "  +6.24%  "
